$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 296.77777
$ws.Range("I11").Value = 296.77777
$ws.Range("K11").Value = 296.77777
$ws.Range("M11").Value = -156.77777
$ws.Range("H86").Value = 9952.916999999999
$ws.Range("I86").Value = 2683
$ws.Range("J86").Value = 12376.223
$ws.Range("K86").Value = 2683
$ws.Range("L86").Value = 12376.223
$ws.Range("M86").Value = -1560
$ws.Range("N86").Value = -14622.223
$ws.Range("H89").Value = 9952.916999999999
$ws.Range("I89").Value = 2683
$ws.Range("J89").Value = 12376.223
$ws.Range("K89").Value = 13415
$ws.Range("L89").Value = 61881.115
$ws.Range("M89").Value = -7799
$ws.Range("N89").Value = -73113.11499999999
$ws.Range("H100").Value = 2676.5386
$ws.Range("I100").Value = 2336.875
$ws.Range("J100").Value = 3220
$ws.Range("K100").Value = 2336.875
$ws.Range("L100").Value = 3220
$ws.Range("M100").Value = -1795.875
$ws.Range("N100").Value = -4302
$ws.Range("H106").Value = 1801.091
$ws.Range("I106").Value = 1325.25
$ws.Range("K106").Value = 1325.25
$ws.Range("M106").Value = -694.25
$ws.Range("H129").Value = 776.9091
$ws.Range("I129").Value = 625.2857
$ws.Range("J129").Value = 847.6667
$ws.Range("K129").Value = 1875.8571
$ws.Range("L129").Value = 2543.0001
$ws.Range("M129").Value = 3124.1429
$ws.Range("N129").Value = -12543.0001
$ws.Range("H138").Value = 2161.6052
$ws.Range("I138").Value = 1334.92
$ws.Range("J138").Value = 3751.3845
$ws.Range("K138").Value = 4004.76
$ws.Range("L138").Value = 11254.1535
$ws.Range("M138").Value = 1135.24
$ws.Range("N138").Value = -21534.1535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18735.457
$ws.Range("I32").Value = 21173.53
$ws.Range("J32").Value = 3192.75
$ws.Range("K32").Value = 21173.53
$ws.Range("L32").Value = 3192.75
$ws.Range("M32").Value = -20886.53
$ws.Range("N32").Value = -3766.75
$ws.Range("H61").Value = 5125
$ws.Range("I61").Value = 2500
$ws.Range("K61").Value = 2500
$ws.Range("M61").Value = -2288
$ws.Range("H97").Value = 1208.3182
$ws.Range("I97").Value = 967.3889
$ws.Range("J97").Value = 2292.5
$ws.Range("K97").Value = 967.3889
$ws.Range("L97").Value = 2292.5
$ws.Range("M97").Value = -471.3889
$ws.Range("N97").Value = -3284.5
$ws.Range("H122").Value = 1369.5
$ws.Range("I122").Value = 1427.1538
$ws.Range("J122").Value = 1219.6
$ws.Range("K122").Value = 4281.4614
$ws.Range("L122").Value = 3658.8
$ws.Range("M122").Value = -1831.4614
$ws.Range("N122").Value = -8558.799999999999
$ws.Range("H136").Value = 5125
$ws.Range("I136").Value = 2500
$ws.Range("K136").Value = 7500
$ws.Range("M136").Value = -4950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2299.7144
$ws.Range("I99").Value = 1813.375
$ws.Range("J99").Value = 2948.1667
$ws.Range("K99").Value = 1813.375
$ws.Range("L99").Value = 2948.1667
$ws.Range("M99").Value = -315.375
$ws.Range("N99").Value = -5944.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11670.974
$ws.Range("I31").Value = 12918.844
$ws.Range("J31").Value = 5015.6665
$ws.Range("K31").Value = 12918.844
$ws.Range("L31").Value = 5015.6665
$ws.Range("M31").Value = -12623.844
$ws.Range("N31").Value = -5605.6665
$ws.Range("H34").Value = 11670.974
$ws.Range("I34").Value = 12918.844
$ws.Range("J34").Value = 5015.6665
$ws.Range("K34").Value = 12918.844
$ws.Range("L34").Value = 5015.6665
$ws.Range("M34").Value = -12716.844
$ws.Range("N34").Value = -5419.6665
$ws.Range("H132").Value = 23308.916
$ws.Range("I132").Value = 30557.354
$ws.Range("K132").Value = 91672.06200000001
$ws.Range("M132").Value = -89142.06200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 1503
$ws.Range("J35").Value = 1503
$ws.Range("L35").Value = 4509
$ws.Range("N35").Value = -5085
$ws.Range("H57").Value = 4666.6665
$ws.Range("I57").Value = 2000
$ws.Range("K57").Value = 6000
$ws.Range("M57").Value = -5441
$ws.Range("H121").Value = 9630.727999999999
$ws.Range("J121").Value = 17267.166
$ws.Range("L121").Value = 51801.49800000001
$ws.Range("N121").Value = -54421.49800000001
$ws.Range("H131").Value = 777.77
$ws.Range("J131").Value = 783.1326
$ws.Range("L131").Value = 2349.3978
$ws.Range("N131").Value = -12429.3978

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 25005250
$ws.Range("J52").Value = 25005250
$ws.Range("L52").Value = 25005250
$ws.Range("N52").Value = -25005768
$ws.Range("H107").Value = 1889
$ws.Range("I107").Value = 400
$ws.Range("J107").Value = 2137.1667
$ws.Range("K107").Value = 400
$ws.Range("L107").Value = 2137.1667
$ws.Range("M107").Value = 1520
$ws.Range("N107").Value = -5977.1667
$ws.Range("H141").Value = 42000
$ws.Range("J141").Value = 42000
$ws.Range("L141").Value = 42000
$ws.Range("N141").Value = -52360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 414.35715
$ws.Range("I16").Value = 394.55554
$ws.Range("J16").Value = 450
$ws.Range("K16").Value = 394.55554
$ws.Range("L16").Value = 450
$ws.Range("M16").Value = -224.55554
$ws.Range("N16").Value = -790
$ws.Range("H48").Value = 20000
$ws.Range("J48").Value = 20000
$ws.Range("L48").Value = 20000
$ws.Range("N48").Value = -21322
$ws.Range("H93").Value = 2275.6296
$ws.Range("I93").Value = 2457.7896
$ws.Range("J93").Value = 1843
$ws.Range("K93").Value = 2457.7896
$ws.Range("L93").Value = 1843
$ws.Range("M93").Value = -1209.7896
$ws.Range("N93").Value = -4339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1613.7
$ws.Range("I107").Value = 694
$ws.Range("K107").Value = 2082
$ws.Range("M107").Value = -162
$ws.Range("H113").Value = 2458082.2
$ws.Range("I113").Value = 1264.2222
$ws.Range("K113").Value = 3792.6666
$ws.Range("M113").Value = -1622.6666
$ws.Range("H126").Value = 1093.6666
$ws.Range("J126").Value = 2083.4285
$ws.Range("L126").Value = 6250.2855
$ws.Range("N126").Value = -11190.2855
$ws.Range("H132").Value = 2410.7856
$ws.Range("I132").Value = 2139.8948
$ws.Range("J132").Value = 2982.6667
$ws.Range("K132").Value = 6419.6844
$ws.Range("L132").Value = 8948.000100000001
$ws.Range("M132").Value = -3889.6844
$ws.Range("N132").Value = -14008.0001
$ws.Range("H135").Value = 44286.25
$ws.Range("I135").Value = 25000
$ws.Range("K135").Value = 25000
$ws.Range("M135").Value = -19930
$ws.Range("H140").Value = 47959.5
$ws.Range("J140").Value = 47959.5
$ws.Range("L140").Value = 47959.5
$ws.Range("N140").Value = -58319.5
$ws.Range("H141").Value = 59942.8
$ws.Range("J141").Value = 59942.8
$ws.Range("L141").Value = 59942.8
$ws.Range("N141").Value = -70302.8
